$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add four new rows to the "Snippets" table (rows 170-173), expanding the
# table / autofilter / dimension ranges automatically.
1..4 | ForEach-Object { $lo.ListRows.Add() | Out-Null }

# Fill in column D (MethodNameInTheSnippet) first ...
$ws.Range("D170").Value = "connectStraightLine"
$ws.Range("D171").Value = "connectStraightLine"
$ws.Range("D172").Value = "disconnectStraightLine"
$ws.Range("D173").Value = "disconnectStraightLine"

# ... then column A (Class) ...
$ws.Range("A170").Value = "Line"
$ws.Range("A171").Value = "Line"
$ws.Range("A172").Value = "Line"
$ws.Range("A173").Value = "Line"

# ... then column B (Method/Prop/Rel Name) ...
$ws.Range("B170").Value = "connectBeginShape"
$ws.Range("B171").Value = "connectEndShape"
$ws.Range("B172").Value = "disconnectBeginShape"
$ws.Range("B173").Value = "disconnectEndShape"

# ... then column C (SnippetIdIntheYAMLFile), which reuses an existing value.
$ws.Range("C170").Value = "excel-shape-lines"
$ws.Range("C171").Value = "excel-shape-lines"
$ws.Range("C172").Value = "excel-shape-lines"
$ws.Range("C173").Value = "excel-shape-lines"

# Columns A, B and D on the new rows pick up an explicit "General" number
# format (matching the new cellXfs style the workbook gained), while column
# C is left with the default style.
$ws.Range("A170:B173").NumberFormat = "General"
$ws.Range("D170:D173").NumberFormat = "General"

# Move the active selection to match where the editor ended up.
$ws.Range("B173").Select() | Out-Null
